# Added 4wk low sales check: re-ran forecast model, which revised the
# near-term (first ~8 weeks) forecast downward after flagging low recent
# sales, and refreshed the dependent inventory coverage / seasonality index
# / reorder-urgency / summary totals that derive from the forecast.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2 (W10)
$ws1.Range("D2").Value = 24
$ws1.Range("H2").Value = 8.699999999999999
$ws1.Range("L2").Value = 1.11

# Row 3 (W11)
$ws1.Range("D3").Value = 24
$ws1.Range("H3").Value = 7.64
$ws1.Range("L3").Value = 1.08

# Row 4 (W12)
$ws1.Range("H4").Value = 6.59
$ws1.Range("L4").Value = 0.83

# Row 5 (W13)
$ws1.Range("D5").Value = 25
$ws1.Range("H5").Value = 5.5
$ws1.Range("L5").Value = 1.01

# Row 6 (W14)
$ws1.Range("D6").Value = 25
$ws1.Range("H6").Value = 4.47
$ws1.Range("L6").Value = 0.96

# Row 7 (W15)
$ws1.Range("D7").Value = 25
$ws1.Range("H7").Value = 3.44
$ws1.Range("L7").Value = 1.1

# Row 8 (W16)
$ws1.Range("D8").Value = 26
$ws1.Range("H8").Value = 2.42
$ws1.Range("L8").Value = 1.09

# Row 9 (W17)
$ws1.Range("D9").Value = 26
$ws1.Range("H9").Value = 1.41
$ws1.Range("J9").Value = "Normal"
$ws1.Range("L9").Value = 0.89

# Row 10 (W18)
$ws1.Range("D10").Value = 26
$ws1.Range("H10").Value = 0.41
$ws1.Range("L10").Value = 0.86

# Row 11 (W19)
$ws1.Range("D11").Value = 26
$ws1.Range("L11").Value = 1.06

# Row 12 (W20)
$ws1.Range("D12").Value = 26
$ws1.Range("L12").Value = 0.95

# Row 13 (W21)
$ws1.Range("D13").Value = 27
$ws1.Range("L13").Value = 0.86

# Row 14 (W22)
$ws1.Range("D14").Value = 27
$ws1.Range("L14").Value = 0.9399999999999999

# Row 15 (W23)
$ws1.Range("L15").Value = 1.1

# Row 16 (W24)
$ws1.Range("D16").Value = 27
$ws1.Range("L16").Value = 1.15

# Row 17 (W25)
$ws1.Range("L17").Value = 1.15

# --- Sheet 2: "Summary" ---
# These cells hold numbers-as-text (t="inlineStr" in the source file), so
# force Text format first to stop Excel from auto-converting the numeric
# string back into a real number on assignment.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "422"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "204"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "100"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "28"
